$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ALBIRENA GARCIA ANGEELO ALONSO", 158),
    @("MANUEL LEUNARDO PRADO BAILON", 147),
    @("URRIOLA ARISMENDIZ INGRID MARYURI", 138),
    @("ALAMA NIMA CLARITZA MABEL", 136),
    @("CORDOVA CARMEN ANGIE NATALLY", 134),
    @("MARYURI OJEDA VALLE", 127),
    @("ATOCHE PALACIOS LUIS ANGEL", 125),
    @("AGURTO ORDINOLA LISBET JAQUELIN", 119),
    @("JUAREZ CARMEN PIERRE ALEXANDER", 111),
    @("VEGA ROBLEDO FERNANDO ERNESTO", 111),
    @("RUIDIAS FRIAS MELISSA VICTORIA", 104),
    @("CARREÑO PALACIOS KATHERINE DE LOS MILAGROS", 100),
    @("ROMAN GALECIO MARITZA DEL PILAR", 100),
    @("BERNAOLA CARMEN ZUMIKO YASHURY", 94)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
